$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1688").Value = 48708
$ws.Range("B1688").Value = 200123
$ws.Range("C1688").Value = "koelkasten en diepvriezers"
$ws.Range("A1689").Value = 48711
$ws.Range("B1689").Value = 200123
$ws.Range("C1689").Value = "afgedankte apparatuur die chloorfluorkoolstoffen bevat"
$ws.Range("A1690").Value = 48714
$ws.Range("B1690").Value = 200123
$ws.Range("C1690").Value = "afgedankte apparatuur die chloorfluorkoolwaterstoffen bevat"
$ws.Range("A1691").Value = 48722
$ws.Range("B1691").Value = 200123
$ws.Range("C1691").Value = "koelkasten"
$ws.Range("A1692").Value = 48737
$ws.Range("B1692").Value = 200123
$ws.Range("C1692").Value = "koelkasten, niet route"
$ws.Range("A1693").Value = 48756
$ws.Range("B1693").Value = 200123
$ws.Range("C1693").Value = "koelkasten per stuk, niet route"
$ws.Range("A1694").Value = 48761
$ws.Range("B1694").Value = 200123
$ws.Range("C1694").Value = "witgoed, niet route"
$ws.Range("A1695").Value = 48762
$ws.Range("B1695").Value = 200123
$ws.Range("C1695").Value = "wit- en bruingoed"
$ws.Range("A1696").Value = 48771
$ws.Range("B1696").Value = 200123
$ws.Range("C1696").Value = "koelkasten/diepvriezers"
$ws.Range("A1697").Value = 48779
$ws.Range("B1697").Value = 200123
$ws.Range("C1697").Value = "witgoed"
$ws.Range("A1698").Value = 48817
$ws.Range("B1698").Value = 200123
$ws.Range("C1698").Value = "koelkasten (witgoed)"
$ws.Range("A1699").Value = 48834
$ws.Range("B1699").Value = 200123
$ws.Range("C1699").Value = "koelkasten, professioneel"
$ws.Range("A1700").Value = 48853
$ws.Range("B1700").Value = 200123
$ws.Range("C1700").Value = "koelkasten, -apparatuur professioneel"
$ws.Range("A1701").Value = 48931
$ws.Range("B1701").Value = 200123
$ws.Range("C1701").Value = "1120 koel- en vries apparatuur"
$ws.Range("A1702").Value = 48939
$ws.Range("B1702").Value = 200123
$ws.Range("C1702").Value = "koelingen"
$ws.Range("A1703").Value = 48952
$ws.Range("B1703").Value = 200123
$ws.Range("C1703").Value = "koel en vries app"
$ws.Range("A1704").Value = 48961
$ws.Range("B1704").Value = 200123
$ws.Range("C1704").Value = "koel en vries apparatuur"
$ws.Range("A1705").Value = 48977
$ws.Range("B1705").Value = 200123
$ws.Range("C1705").Value = "koel en vries app, waterkoelers"
$ws.Range("A1706").Value = 48980
$ws.Range("B1706").Value = 200123
$ws.Range("C1706").Value = "koel en vries apparaten"
$ws.Range("A1707").Value = 48999
$ws.Range("B1707").Value = 200123
$ws.Range("C1707").Value = "radarsnelheidapparatuur"
$ws.Range("A1708").Value = 49000
$ws.Range("B1708").Value = 200123
$ws.Range("C1708").Value = "inbouwapparatuur radarvoertuigen"
$ws.Range("A1709").Value = 49004
$ws.Range("B1709").Value = 200123
$ws.Range("C1709").Value = "koel en vrieskisten"
$ws.Range("A1710").Value = 49035
$ws.Range("B1710").Value = 200123
$ws.Range("C1710").Value = "koelkasten/diepvriezers/koelinstallatie"
$ws.Range("A1711").Value = 49037
$ws.Range("B1711").Value = 200123
$ws.Range("C1711").Value = "airco unit(s)"
$ws.Range("A1712").Value = 49041
$ws.Range("B1712").Value = 200123
$ws.Range("C1712").Value = "witgoed (cfk) - bruingoed - electronica"
$ws.Range("A1713").Value = 49154
$ws.Range("B1713").Value = 200134
$ws.Range("C1713").Value = "niet onder 20 01 33 vallende batterijen en accus"
$ws.Range("A1714").Value = 49155
$ws.Range("B1714").Value = 200134
$ws.Range("C1714").Value = "niet onder 20 01 33 vallende batterijen en accu's"
$ws.Range("A1715").Value = 49185
$ws.Range("B1715").Value = 200110
$ws.Range("C1715").Value = "kleding"
$ws.Range("A1716").Value = 49187
$ws.Range("B1716").Value = 200110
$ws.Range("C1716").Value = "bedrijfsafval cat. c"
$ws.Range("A1717").Value = 49188
$ws.Range("B1717").Value = 200110
$ws.Range("C1717").Value = "textiel"
$ws.Range("A1718").Value = 49189
$ws.Range("B1718").Value = 200110
$ws.Range("C1718").Value = "textiel, ter vernietiging"
$ws.Range("A1719").Value = 49190
$ws.Range("B1719").Value = 40109
$ws.Range("C1719").Value = "afval van bewerking en afwerking"
$ws.Range("A1720").Value = 49191
$ws.Range("B1720").Value = 40109
$ws.Range("C1720").Value = "melamine-ureum/formaline co polymeer"
$ws.Range("A1721").Value = 49192
$ws.Range("B1721").Value = 40222
$ws.Range("C1721").Value = "scrap"
$ws.Range("A1722").Value = 49193
$ws.Range("B1722").Value = 40209
$ws.Range("C1722").Value = "afval van composietmaterialen (ge?mpregneerde textiel, elastomeren, plastomeren)"
$ws.Range("A1723").Value = 49198
$ws.Range("B1723").Value = 40209
$ws.Range("C1723").Value = "afval van composietmaterialen (geimpregneerde textiel, elastomeren, plastomeren)"
$ws.Range("A1724").Value = 49199
$ws.Range("B1724").Value = 40209
$ws.Range("C1724").Value = "linoleum met pvc"
$ws.Range("A1725").Value = 49200
$ws.Range("B1725").Value = 40220
$ws.Range("C1725").Value = "slib van productie van polyprop"
$ws.Range("A1726").Value = 49206
$ws.Range("B1726").Value = 40220
$ws.Range("C1726").Value = "industrieel slib (ba)"
$ws.Range("A1727").Value = 49219
$ws.Range("B1727").Value = 40220
$ws.Range("C1727").Value = "zuiveringslib van tapijtindustrie"
$ws.Range("A1728").Value = 49229
$ws.Range("B1728").Value = 40220
$ws.Range("C1728").Value = "slib"
$ws.Range("A1729").Value = 49230
$ws.Range("B1729").Value = 40106
$ws.Range("C1729").Value = "industrieel slib (ba)"
$ws.Range("A1730").Value = 49231
$ws.Range("B1730").Value = 40102
$ws.Range("C1730").Value = "loogafval"
$ws.Range("A1731").Value = 49235
$ws.Range("B1731").Value = 40216
$ws.Range("C1731").Value = "kleurstoffen en pigmenten die gevaarlijke stoffen bevatten"
$ws.Range("A1732").Value = 49257
$ws.Range("B1732").Value = 40219
$ws.Range("C1732").Value = "slib van afvalwaterbehandeling ter plaatse dat gevaarlijke stoffen bevat"
$ws.Range("A1733").Value = 49260
$ws.Range("B1733").Value = 40214
$ws.Range("C1733").Value = "afval van afwerking dat organische oplosmiddelen bevat"

$ws.Rows("1734:1752").Delete()
